# Replace the inline "Minimum site area for Left-Behind Plots" picture
# with a hyperlink run pointing at the same image, now hosted externally
# on ura.gov.sg (styled with the built-in "Hyperlink" character style).

$d = $word.ActiveDocument

$url = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Flats-Condominiums/F02_Minimum_Site_Area.jpg?h=100%25&w=100%25"

# Locate the inline picture by its alt text (there is exactly one in this
# document) rather than assuming it is InlineShapes(1), so the script keeps
# working even if unrelated shapes are added earlier in the body.
$target = $null
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $candidate = $d.InlineShapes.Item($i)
    if ($candidate.AlternativeText -eq "Minimum site area for Left-Behind Plots") {
        $target = $candidate
    }
}

if ($target -ne $null) {
    $pictureRange = $target.Range
    $target.Delete()
    $d.Hyperlinks.Add($pictureRange, $url, $null, $null, $url) | Out-Null
}

Write-Output "done"
